# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table switches from the built-in table style
#    {1E065C92-22BC-4D9A-91B0-BF367D49244A} to
#    {7CCC535D-9383-4E0C-850A-D0281F5A88E8}.
#
# 2) The deck's applied theme ("Integral" / Red Violet) is swapped for the
#    theme that used to belong only to the Notes Master ("Office Theme" /
#    Office blue) - i.e. the live theme's 12 scheme colors are changed to
#    the classic Office palette.

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$tableShape = $null
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    }
}
$tableShape.Table.ApplyStyle("{7CCC535D-9383-4E0C-850A-D0281F5A88E8}")

# --- 2. Swap the live theme's colour scheme for the Office Theme palette --
$officeColors = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = HexToRgb($officeColors[$i])
}
